$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.363865666666667
$ws.Range("H2").Value = 4.091597
$ws.Range("I2").Value = 0.35258381842799
$ws.Range("J2").Value = 0.35258381842799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 18.382477
$ws.Range("N2").Value = 55.147431
$ws.Range("O2").Value = 0.06380158579420245
$ws.Range("P2").Value = 0.06380158579420243
$ws.Range("Q2").Value = 25.07122924858966
$ws.Range("R2").Value = 225.641063237307
$ws.Range("S2").Value = 0.0224954067410809
$ws.Range("T2").Value = 0.02249540674108089

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.363865666666667
$ws.Range("H3").Value = 4.091597
$ws.Range("I3").Value = 0.35258381842799
$ws.Range("J3").Value = 0.35258381842799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 82.64333833333332
$ws.Range("N3").Value = 247.930015
$ws.Range("O3").Value = 0.2868370808239535
$ws.Range("P3").Value = 0.2868370808239535
$ws.Range("Q3").Value = 112.7144117315505
$ws.Range("R3").Value = 1014.429705583955
$ws.Range("S3").Value = 0.1011341132236475
$ws.Range("T3").Value = 0.1011341132236475

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.363865666666667
$ws.Range("H4").Value = 4.091597
$ws.Range("I4").Value = 0.35258381842799
$ws.Range("J4").Value = 0.35258381842799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 91.25099666666667
$ws.Range("N4").Value = 273.75299
$ws.Range("O4").Value = 0.3167123936907314
$ws.Range("P4").Value = 0.3167123936907314
$ws.Range("Q4").Value = 124.4541014027811
$ws.Range("R4").Value = 1120.08691262503
$ws.Range("S4").Value = 0.1116676651109469
$ws.Range("T4").Value = 0.1116676651109469

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.363865666666667
$ws.Range("H5").Value = 4.091597
$ws.Range("I5").Value = 0.35258381842799
$ws.Range("J5").Value = 0.35258381842799
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 81.28845566666666
$ws.Range("N5").Value = 243.865367
$ws.Range("O5").Value = 0.2821345773094157
$ws.Range("P5").Value = 0.2821345773094157
$ws.Range("Q5").Value = 110.8665337801221
$ws.Range("R5").Value = 997.798804021099
$ws.Range("S5").Value = 0.09947608657832072
$ws.Range("T5").Value = 0.09947608657832072

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.363865666666667
$ws.Range("H6").Value = 4.091597
$ws.Range("I6").Value = 0.35258381842799
$ws.Range("J6").Value = 0.35258381842799
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 14.55416966666667
$ws.Range("N6").Value = 43.662509
$ws.Range("O6").Value = 0.0505143623816971
$ws.Range("P6").Value = 0.0505143623816971
$ws.Range("Q6").Value = 19.84993231520811
$ws.Range("R6").Value = 178.649390836873
$ws.Range("S6").Value = 0.01781054677399398
$ws.Range("T6").Value = 0.01781054677399398

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.247734666666667
$ws.Range("H7").Value = 3.743204
$ws.Range("I7").Value = 0.3225618650798028
$ws.Range("J7").Value = 0.3225618650798028
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.382477
$ws.Range("N7").Value = 55.147431
$ws.Range("O7").Value = 0.06380158579420245
$ws.Range("P7").Value = 0.06380158579420243
$ws.Range("Q7").Value = 22.93645381210267
$ws.Range("R7").Value = 206.428084308924
$ws.Range("S7").Value = 0.020579958508827
$ws.Range("T7").Value = 0.02057995850882699

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.247734666666667
$ws.Range("H8").Value = 3.743204
$ws.Range("I8").Value = 0.3225618650798028
$ws.Range("J8").Value = 0.3225618650798028
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 82.64333833333332
$ws.Range("N8").Value = 247.930015
$ws.Range("O8").Value = 0.2868370808239535
$ws.Range("P8").Value = 0.2868370808239535
$ws.Range("Q8").Value = 103.1169582075622
$ws.Range("R8").Value = 928.05262386806
$ws.Range("S8").Value = 0.09252270376462059
$ws.Range("T8").Value = 0.09252270376462057

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.247734666666667
$ws.Range("H9").Value = 3.743204
$ws.Range("I9").Value = 0.3225618650798028
$ws.Range("J9").Value = 0.3225618650798028
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 91.25099666666667
$ws.Range("N9").Value = 273.75299
$ws.Range("O9").Value = 0.3167123936907314
$ws.Range("P9").Value = 0.3167123936907314
$ws.Range("Q9").Value = 113.8570319088845
$ws.Range("R9").Value = 1024.71328717996
$ws.Range("S9").Value = 0.1021593404027711
$ws.Range("T9").Value = 0.1021593404027711

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.247734666666667
$ws.Range("H10").Value = 3.743204
$ws.Range("I10").Value = 0.3225618650798028
$ws.Range("J10").Value = 0.3225618650798028
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 81.28845566666666
$ws.Range("N10").Value = 243.865367
$ws.Range("O10").Value = 0.2821345773094157
$ws.Range("P10").Value = 0.2821345773094157
$ws.Range("Q10").Value = 101.4264241350965
$ws.Range("R10").Value = 912.8378172158681
$ws.Range("S10").Value = 0.09100585546042693
$ws.Range("T10").Value = 0.09100585546042692

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.247734666666667
$ws.Range("H11").Value = 3.743204
$ws.Range("I11").Value = 0.3225618650798028
$ws.Range("J11").Value = 0.3225618650798028
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 14.55416966666667
$ws.Range("N11").Value = 43.662509
$ws.Range("O11").Value = 0.0505143623816971
$ws.Range("P11").Value = 0.0505143623816971
$ws.Range("Q11").Value = 18.15974203764845
$ws.Range("R11").Value = 163.437678338836
$ws.Range("S11").Value = 0.01629400694315725
$ws.Range("T11").Value = 0.01629400694315724

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.808894
$ws.Range("H12").Value = 2.426682
$ws.Range("I12").Value = 0.2091136555409713
$ws.Range("J12").Value = 0.2091136555409713
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.382477
$ws.Range("N12").Value = 55.147431
$ws.Range("O12").Value = 0.06380158579420245
$ws.Range("P12").Value = 0.06380158579420243
$ws.Range("Q12").Value = 14.869475350438
$ws.Range("R12").Value = 133.825278153942
$ws.Range("S12").Value = 0.01334178283473658
$ws.Range("T12").Value = 0.01334178283473658

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.808894
$ws.Range("H13").Value = 2.426682
$ws.Range("I13").Value = 0.2091136555409713
$ws.Range("J13").Value = 0.2091136555409713
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 82.64333833333332
$ws.Range("N13").Value = 247.930015
$ws.Range("O13").Value = 0.2868370808239535
$ws.Range("P13").Value = 0.2868370808239535
$ws.Range("Q13").Value = 66.84970051780333
$ws.Range("R13").Value = 601.6473046602299
$ws.Range("S13").Value = 0.05998155051579795
$ws.Range("T13").Value = 0.05998155051579795

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.808894
$ws.Range("H14").Value = 2.426682
$ws.Range("I14").Value = 0.2091136555409713
$ws.Range("J14").Value = 0.2091136555409713
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 91.25099666666667
$ws.Range("N14").Value = 273.75299
$ws.Range("O14").Value = 0.3167123936907314
$ws.Range("P14").Value = 0.3167123936907314
$ws.Range("Q14").Value = 73.81238369768667
$ws.Range("R14").Value = 664.31145327918
$ws.Range("S14").Value = 0.0662288863998001
$ws.Range("T14").Value = 0.0662288863998001

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.808894
$ws.Range("H15").Value = 2.426682
$ws.Range("I15").Value = 0.2091136555409713
$ws.Range("J15").Value = 0.2091136555409713
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 81.28845566666666
$ws.Range("N15").Value = 243.865367
$ws.Range("O15").Value = 0.2821345773094157
$ws.Range("P15").Value = 0.2821345773094157
$ws.Range("Q15").Value = 65.75374405803267
$ws.Range("R15").Value = 591.783696522294
$ws.Range("S15").Value = 0.05899819281567868
$ws.Range("T15").Value = 0.05899819281567868

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.808894
$ws.Range("H16").Value = 2.426682
$ws.Range("I16").Value = 0.2091136555409713
$ws.Range("J16").Value = 0.2091136555409713
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 14.55416966666667
$ws.Range("N16").Value = 43.662509
$ws.Range("O16").Value = 0.0505143623816971
$ws.Range("P16").Value = 0.0505143623816971
$ws.Range("Q16").Value = 11.77278051834867
$ws.Range("R16").Value = 105.955024665138
$ws.Range("S16").Value = 0.01056324297495801
$ws.Range("T16").Value = 0.01056324297495801

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.4477083333333334
$ws.Range("H17").Value = 1.343125
$ws.Range("I17").Value = 0.1157406609512359
$ws.Range("J17").Value = 0.1157406609512359
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 18.382477
$ws.Range("N17").Value = 55.147431
$ws.Range("O17").Value = 0.06380158579420245
$ws.Range("P17").Value = 0.06380158579420243
$ws.Range("Q17").Value = 8.229988140208333
$ws.Range("R17").Value = 74.06989326187501
$ws.Range("S17").Value = 0.007384437709557976
$ws.Range("T17").Value = 0.007384437709557975

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 0.4477083333333334
$ws.Range("H18").Value = 1.343125
$ws.Range("I18").Value = 0.1157406609512359
$ws.Range("J18").Value = 0.1157406609512359
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 82.64333833333332
$ws.Range("N18").Value = 247.930015
$ws.Range("O18").Value = 0.2868370808239535
$ws.Range("P18").Value = 0.2868370808239535
$ws.Range("Q18").Value = 37.00011126631944
$ws.Range("R18").Value = 333.001001396875
$ws.Range("S18").Value = 0.03319871331988745
$ws.Range("T18").Value = 0.03319871331988745

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 0.4477083333333334
$ws.Range("H19").Value = 1.343125
$ws.Range("I19").Value = 0.1157406609512359
$ws.Range("J19").Value = 0.1157406609512359
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 91.25099666666667
$ws.Range("N19").Value = 273.75299
$ws.Range("O19").Value = 0.3167123936907314
$ws.Range("P19").Value = 0.3167123936907314
$ws.Range("Q19").Value = 40.85383163263889
$ws.Range("R19").Value = 367.6844846937501
$ws.Range("S19").Value = 0.03665650177721329
$ws.Range("T19").Value = 0.03665650177721329

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 0.4477083333333334
$ws.Range("H20").Value = 1.343125
$ws.Range("I20").Value = 0.1157406609512359
$ws.Range("J20").Value = 0.1157406609512359
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 81.28845566666666
$ws.Range("N20").Value = 243.865367
$ws.Range("O20").Value = 0.2821345773094157
$ws.Range("P20").Value = 0.2821345773094157
$ws.Range("Q20").Value = 36.39351900576389
$ws.Range("R20").Value = 327.541671051875
$ws.Range("S20").Value = 0.03265444245498934
$ws.Range("T20").Value = 0.03265444245498934

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 0.4477083333333334
$ws.Range("H21").Value = 1.343125
$ws.Range("I21").Value = 0.1157406609512359
$ws.Range("J21").Value = 0.1157406609512359
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 14.55416966666667
$ws.Range("N21").Value = 43.662509
$ws.Range("O21").Value = 0.0505143623816971
$ws.Range("P21").Value = 0.0505143623816971
$ws.Range("Q21").Value = 6.516023044513889
$ws.Range("R21").Value = 6.516023044513889
$ws.Range("S21").Value = 0.005846565689587872
$ws.Range("T21").Value = 0.005846565689587871
